$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 contents (Case 4 log entry)
$ws.Range("A5").Value = "Satuday 12:00PM - 2:00 PM"
$ws.Range("B5").Value = "Understanding the program/ Looking over case 4"
$ws.Range("C5").Value = "Fixed offset on Lw and sw instructions. Fixing other potential issues"
$ws.Range("D5").Value = "Rami"

# Widen column C to fit the new, longer progress note
$ws.Columns.Item(3).ColumnWidth = 59.5
